# Auto-generated Excel COM-interop script applying the cryptos.xlsx price/volume update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''27.001.28'
$ws.Range('E2').Value = '  -1.35%  '
$ws.Range('D3').Value = '''1.820.21'
$ws.Range('E3').Value = '  -1.08%  '
$ws.Range('D4').Value = '''1.010'
$ws.Range('E4').Value = '  -0.48%  '
$ws.Range('D5').Value = '''309.93'
$ws.Range('E5').Value = '  -1.60%  '
$ws.Range('E6').Value = '  -0.43%  '
$ws.Range('D7').Value = '''0.4636'
$ws.Range('E7').Value = '  -2.35%  '
$ws.Range('D8').Value = '''0.3638'
$ws.Range('E8').Value = '  -1.71%  '
$ws.Range('D9').Value = '''0.07291'
$ws.Range('E9').Value = '  -2.31%  '
$ws.Range('D10').Value = '''0.8657'
$ws.Range('E10').Value = '  -2.25%  '
$ws.Range('D11').Value = '''19.82'
$ws.Range('E11').Value = '  -3.32%  '
$ws.Range('D12').Value = '''1.879.09'
$ws.Range('E12').Value = '  +1.61%  '
$ws.Range('D13').Value = '''0.07572'
$ws.Range('E13').Value = '  +2.70%  '
$ws.Range('D14').Value = '''93.10'
$ws.Range('E14').Value = '  -0.21%  '
$ws.Range('E15').Value = '  -2.61%  '
$ws.Range('D16').Value = '''6.465'
$ws.Range('E16').Value = '  -1.87%  '
$ws.Range('D17').Value = '''1.009'
$ws.Range('E17').Value = '  -0.66%  '
$ws.Range('D18').Value = '''0.000008637'
$ws.Range('E18').Value = '  -2.45%  '
$ws.Range('D19').Value = '''1.009'
$ws.Range('E19').Value = '  -0.43%  '
$ws.Range('D20').Value = '''27.367.14'
$ws.Range('E20').Value = '  -0.10%  '
$ws.Range('D21').Value = '''14.48'
$ws.Range('E21').Value = '  -2.48%  '
$ws.Range('D22').Value = '''5.157'
$ws.Range('E22').Value = '  -3.66%  '
$ws.Range('E23').Value = '  -1.50%  '
$ws.Range('D24').Value = '''2.101.16'
$ws.Range('E24').Value = '  +1.31%  '
$ws.Range('D25').Value = '''151.83'
$ws.Range('E25').Value = '  -0.36%  '
$ws.Range('D26').Value = '''1.866'
$ws.Range('E26').Value = '  -2.39%  '
$ws.Range('D27').Value = '''18.16'
$ws.Range('E27').Value = '  -2.62%  '
$ws.Range('E28').Value = '  -3.38%  '
$ws.Range('D31').Value = '''0.08902'
$ws.Range('E31').Value = '  -0.82%  '
$ws.Range('D32').Value = '''2.953'
$ws.Range('E32').Value = '  +0.36%  '
$ws.Range('D33').Value = '''0.7305'
$ws.Range('E33').Value = '  -4.02%  '
$ws.Range('D34').Value = '''1.139'
$ws.Range('E34').Value = '  -3.31%  '
$ws.Range('E35').Value = '  -3.00%  '
$ws.Range('D36').Value = '''1.010'
$ws.Range('E36').Value = '  -0.40%  '
$ws.Range('D37').Value = '''2.517'
$ws.Range('E37').Value = '  +5.61%  '
$ws.Range('D38').Value = '''1.073'
$ws.Range('E38').Value = '  -3.03%  '
$ws.Range('D39').Value = '''0.05263'
$ws.Range('E39').Value = '  -2.23%  '
$ws.Range('D40').Value = '''0.01919'
$ws.Range('E40').Value = '  -2.29%  '
$ws.Range('D41').Value = '''2.926'
$ws.Range('E41').Value = '  -2.47%  '
$ws.Range('D42').Value = '''7.129'
$ws.Range('E42').Value = '  -2.53%  '
$ws.Range('D43').Value = '''0.5207'
$ws.Range('E43').Value = '  -2.84%  '
$ws.Range('D44').Value = '''0.1632'
$ws.Range('E44').Value = '  -2.11%  '
$ws.Range('D45').Value = '''8.229'
$ws.Range('E45').Value = '  -3.88%  '
$ws.Range('E46').Value = '  -2.79%  '
$ws.Range('D47').Value = '''1.009'
$ws.Range('E47').Value = '  -0.52%  '
$ws.Range('D48').Value = '''10.09'
$ws.Range('E48').Value = '  -4.55%  '
$ws.Range('E49').Value = '  -1.89%  '
$ws.Range('D50').Value = '''1.635'
$ws.Range('E50').Value = '  -2.81%  '
$ws.Range('D51').Value = '''0.06221'

# Rows 29/30 swapped position (InternetComputer(DFINITY) moved above BitcoinCash)
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').Value = '''5.078'
$ws.Range('E29').Value = '  -3.53%  '

$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').Value = '''115.81'
$ws.Range('E30').Value = '  -1.95%  '
